$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete data rows 2-5 (the four obsolete entries), shifting remaining rows up
$ws.Range("A2:B5").EntireRow.Delete() | Out-Null

# Update the selection to match the new data extent
$ws.Range("A2:A13").Select() | Out-Null
